$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $Value, $DonorRef)
    $Sheet.Range($CellRef).NumberFormat = "@"
    $Sheet.Range($CellRef).Value = $Value
    $Sheet.Range($CellRef).Style = $Sheet.Range($DonorRef).Style
}

$ws.Range("D2").Value = "65.856.86"
$ws.Range("E2").Value = "  +6.38%  "
$ws.Range("D3").Value = "3.008.39"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue $ws "D5" "582.79" "B5"
$ws.Range("E5").Value = "  +2.37%  "
Set-TextValue $ws "D6" "162.73" "B6"
$ws.Range("E6").Value = "  +13.24%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").Value = "3.004.80"
$ws.Range("E9").Value = "  +3.58%  "
Set-TextValue $ws "D10" "6.82" "B10"
$ws.Range("E10").Value = "  -1.83%  "
Set-TextValue $ws "D11" "0.156" "B11"
$ws.Range("E11").Value = "  +7.38%  "
Set-TextValue $ws "D12" "0.460" "B12"
$ws.Range("E12").Value = "  +7.31%  "
Set-TextValue $ws "D13" "0.0000250" "B13"
$ws.Range("E13").Value = "  +8.38%  "
Set-TextValue $ws "D14" "34.73" "B14"
$ws.Range("E14").Value = "  +8.04%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "65.898.08"
$ws.Range("E16").Value = "  +6.55%  "
$ws.Range("D17").Value = "3.509.94"
$ws.Range("E17").Value = "  +3.77%  "
Set-TextValue $ws "D18" "6.96" "B18"
$ws.Range("D19").Value = "3.011.80"
$ws.Range("E19").Value = "  +3.90%  "
Set-TextValue $ws "D20" "457.37" "B20"
$ws.Range("E20").Value = "  +6.41%  "
Set-TextValue $ws "D21" "13.95" "B21"
$ws.Range("E21").Value = "  +8.01%  "
$ws.Range("E22").Value = "  +5.49%  "
Set-TextValue $ws "D23" "7.36" "B23"
$ws.Range("E23").Value = "  +7.30%  "
Set-TextValue $ws "D24" "82.51" "B24"
$ws.Range("E24").Value = "  +4.43%  "
Set-TextValue $ws "D25" "2.30" "B25"
$ws.Range("E25").Value = "  +13.30%  "
Set-TextValue $ws "D26" "12.35" "B26"
$ws.Range("E26").Value = "  +3.17%  "
Set-TextValue $ws "D27" "10.62" "B27"
$ws.Range("E27").Value = "  +5.13%  "
Set-TextValue $ws "D29" "8.18" "B29"
$ws.Range("E29").Value = "  +16.81%  "
$ws.Range("E30").Value = "  +15.68%  "
Set-TextValue $ws "D31" "2.60" "B31"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("E32").Value = "  -6.92%  "
Set-TextValue $ws "D33" "26.99" "B33"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("E34").Value = "  +3.66%  "
Set-TextValue $ws "D36" "0.991" "B36"
$ws.Range("E36").Value = "  +3.77%  "
Set-TextValue $ws "D37" "5.80" "B37"
$ws.Range("E37").Value = "  +7.54%  "
$ws.Range("E38").Value = "  +11.92%  "
Set-TextValue $ws "D39" "3.00" "B39"
$ws.Range("E39").Value = "  +5.71%  "
Set-TextValue $ws "D40" "49.76" "B40"
$ws.Range("E40").Value = "  +1.89%  "
Set-TextValue $ws "D41" "0.308" "B41"
$ws.Range("E41").Value = "  +14.93%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D42" "0.121" "B42"
$ws.Range("E42").Value = "  +6.13%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws "D43" "43.69" "B43"
$ws.Range("E43").Value = "  +7.22%  "
Set-TextValue $ws "D44" "8.47" "B44"
$ws.Range("E44").Value = "  +3.87%  "
Set-TextValue $ws "D45" "387.85" "B45"
$ws.Range("E45").Value = "  +12.03%  "
$ws.Range("D47").Value = "2.796.29"
$ws.Range("E47").Value = "  +3.37%  "
Set-TextValue $ws "D48" "134.91" "B48"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("E49").Value = "  -0.02%  "
Set-TextValue $ws "D50" "23.97" "B50"
$ws.Range("E50").Value = "  +11.18%  "
$ws.Range("E51").Value = "  +3.96%  "
